$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update version number (B5)
$ws.Range("B5").Value = "0002"

# Update Pos-condición text (B19)
$ws.Range("B19").Value = "Se completó el formulario con toda la información correspondiente al actor"

# Let Excel auto-fit the row height now that the text is shorter
$ws.Rows.Item(19).AutoFit() | Out-Null

# Move the active selection to C19, matching the final saved view state
$ws.Range("C19").Select() | Out-Null
